# [SAU-1563] combine ices div and subdiv
#
# The "Catch Data" sheet had two adjacent header columns:
#   H: "ICES division"
#   I: "ICES subdivision"
# They are combined into a single column "ICES area" (kept in column H),
# and the now-redundant "ICES subdivision" column (I) is removed, shifting
# every column after it one place to the left (J..AB -> I..AA).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column I ("ICES subdivision") entirely - Excel shifts J:AB left to I:AA.
$ws.Columns("I").Delete()

# Rename the surviving column header (old H "ICES division") to the merged name.
$ws.Range("H1").Value = "ICES area"
